$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'extra long basketball shorts for men'
$ws.Range("A2").Value = 'compression pants running'
$ws.Range("A3").Value = 'boys baseball compression sleeve'
$ws.Range("A4").Value = 'stretch mark for men'
$ws.Range("A5").Value = 'compression shorts men pack'
$ws.Range("A6").Value = 'knee pads youth girls'
$ws.Range("A7").Value = 'knee sleeves basketball'
$ws.Range("A8").Value = 'xl youth baseball pants'
$ws.Range("A9").Value = 'biking knee sleeve'
$ws.Range("A10").Value = 'boy baseball pants'
$ws.Range("A11").Value = 'men black compression pants'
$ws.Range("A12").Value = 'knee pads exercise'
$ws.Range("A13").Value = 'youth basketball compression shorts'
$ws.Range("A14").Value = 'hockey pads'
$ws.Range("A15").Value = 'calf strain compression sleeve'
$ws.Range("A16").Value = 'mens degree sport'
$ws.Range("A17").Value = 'men long shorts below the knee'
$ws.Range("A18").Value = 'knee pads construction gel'
$ws.Range("A19").Value = 'work knee pads construction'
$ws.Range("A20").Value = 'padded football sleeve'
$ws.Range("A21").Value = 'spandex tights'
$ws.Range("A22").Value = 'adult black baseball pants'
$ws.Range("A23").Value = 'girls baseball pants'
$ws.Range("A24").Value = 'mens small black baseball pants'
$ws.Range("A25").Value = 'youth softball compression sleeve'
$ws.Range("A26").Value = 'padded shorts football'
$ws.Range("A27").Value = '3/4 shorts for men'
$ws.Range("A28").Value = 'silicone strips for stretch marks'
$ws.Range("A29").Value = 'calf compression sleeve basketball'
$ws.Range("A30").Value = 'gel knee pads for construction'
$ws.Range("A31").Value = 'youth baseball short pants'
$ws.Range("A32").Value = 'exercise knee pad'
$ws.Range("A33").Value = 'girls basketball shorts size 6'
$ws.Range("A34").Value = 'mens short tights'
$ws.Range("A35").Value = 'boys knee length shorts'
$ws.Range("A36").Value = 'rodilleras de volleyball'
$ws.Range("A37").Value = 'rodilleras volleyball'
$ws.Range("A38").Value = 'compression running pants men'
$ws.Range("A39").Value = 'knee pads cheap'
$ws.Range("A40").Value = 'exercise knee pads'
$ws.Range("A41").Value = 'pad knee'
$ws.Range("A42").Value = 'compression workout pants men'
$ws.Range("A43").Value = 'compression sleeve for bursitis'
$ws.Range("A44").Value = 'mens running tights shorts'
$ws.Range("A45").Value = 'mens compression girdle'
$ws.Range("A46").Value = 'knee pads outdoor'
$ws.Range("A47").Value = 'motorcycle knee pads for men'
$ws.Range("A48").Value = 'water pants men'
$ws.Range("A49").Value = 'pant baseball men'
$ws.Range("A50").Value = 'volleyball knee pad'
$ws.Range("A51").Value = 'knee pads for volleyball'
$ws.Range("A52").Value = 'pro tights men'
$ws.Range("A53").Value = 'sleeve knee pads'
$ws.Range("A54").Value = 'womens compression leggings'
$ws.Range("A55").Value = 'airsoft knee pads'
$ws.Range("A56").Value = 'skateboard knee pads'
$ws.Range("A57").Value = 'starter youth compression pants'
$ws.Range("A58").Value = 'trolls knee pads'
$ws.Range("A59").Value = 'mens compression pants marvel'
$ws.Range("A60").Value = 'mouthguard basketball youth'
$ws.Range("A61").Value = 'mcdavid youth knee pads'
$ws.Range("A62").Value = 'nike compression leggings'
$ws.Range("A63").Value = 'jordan flight mens basketball pants'
$ws.Range("A64").Value = 'athletic compression pants'
$ws.Range("A65").Value = 'nike pro compression leggings men'
$ws.Range("A66").Value = 'emoji knee pads'
$ws.Range("A67").Value = 'mens compression tights 3 4'
$ws.Range("A68").Value = 'mens compression tights nike'
$ws.Range("A69").Value = 'mens compression tights white'
$ws.Range("A70").Value = 'eastbay compression pants'
$ws.Range("A71").Value = 'women compression leggings'
$ws.Range("A72").Value = 'ladies compression pants'
$ws.Range("A73").Value = 'womans compression leggings'
$ws.Range("A74").Value = 'mcdavid compression pants'
$ws.Range("A75").Value = 'knee brace basketball youth'
$ws.Range("A76").Value = 'basketball kids knee pads'
$ws.Range("A77").Value = 'basketball knee pads mcdavid'
$ws.Range("A78").Value = 'basketball youth jersey'
$ws.Range("A79").Value = 'morris compression knee pads'
$ws.Range("A80").Value = 'wonens compression leggings'
$ws.Range("A81").Value = 'youth baketball knee pads'
$ws.Range("A82").Value = 'olympic mens basketball'
$ws.Range("A83").Value = 'track leggings men'
$ws.Range("A84").Value = 'nike pro dry mens basketball tights'
$ws.Range("A85").Value = 'goalie knee protectors'
$ws.Range("A86").Value = 'basketball clothes for men'
$ws.Range("A87").Value = 'mens running thermal compression pants'
$ws.Range("A88").Value = 'capri pants for men adidas'
$ws.Range("A89").Value = 'men running tights'
$ws.Range("A90").Value = 'mens running tight'
$ws.Range("A91").Value = 'men tights nike'
$ws.Range("A92").Value = 'men tights short'
$ws.Range("A93").Value = 'mens tights dance'
$ws.Range("A94").Value = 'men running tight'
$ws.Range("A95").Value = 'men tights green'
$ws.Range("A96").Value = 'men tights pack'
$ws.Range("A97").Value = 'mens tights grey'
$ws.Range("A98").Value = 'mens tights nike'
$ws.Range("A99").Value = 'mens tights pink'
$ws.Range("A100").Value = 'men legging nike'
